$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.301.31"
$ws.Range("E2").Value = "  -3.39%  "
$ws.Range("D3").Value = "2.479.51"
$ws.Range("E3").Value = "  -2.43%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'313.33"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("D6").Value = "'94.58"
$ws.Range("E6").Value = "  -5.45%  "
$ws.Range("D7").Value = "'0.550"
$ws.Range("E7").Value = "  -2.90%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "'0.500"
$ws.Range("E9").Value = "  -4.29%  "
$ws.Range("D10").Value = "'33.55"
$ws.Range("E10").Value = "  -5.18%  "
$ws.Range("D11").Value = "'0.0783"
$ws.Range("E11").Value = "  -2.73%  "
$ws.Range("E12").Value = "  -0.71%  "
$ws.Range("D13").Value = "'7.01"
$ws.Range("E13").Value = "  -3.96%  "
$ws.Range("D14").Value = "2.860.48"
$ws.Range("E14").Value = "  -2.48%  "
$ws.Range("D15").Value = "'15.38"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("D16").Value = "2.446.62"
$ws.Range("E16").Value = "  -3.84%  "
$ws.Range("E17").Value = "  -2.89%  "
$ws.Range("D18").Value = "41.336.16"
$ws.Range("E18").Value = "  -3.35%  "
$ws.Range("E19").Value = "  -5.91%  "
$ws.Range("D20").Value = "0.0₃0926"
$ws.Range("E20").Value = "  -2.48%  "
$ws.Range("D21").Value = "'11.25"
$ws.Range("E21").Value = "  -8.70%  "
$ws.Range("D22").Value = "'68.76"
$ws.Range("E22").Value = "  -1.43%  "
$ws.Range("D23").Value = "'237.34"
$ws.Range("E23").Value = "  -2.30%  "
$ws.Range("D24").Value = "'2.76"
$ws.Range("E24").Value = "  -3.99%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("E26").Value = "  -5.85%  "
$ws.Range("D27").Value = "'24.11"
$ws.Range("E27").Value = "  -6.16%  "
$ws.Range("D28").Value = "'2.24"
$ws.Range("E28").Value = "  -4.35%  "
$ws.Range("D29").Value = "'9.71"
$ws.Range("E29").Value = "  -4.12%  "
$ws.Range("E30").Value = "  -4.92%  "
$ws.Range("D31").Value = "'152.57"
$ws.Range("E31").Value = "  -3.24%  "
$ws.Range("E32").Value = "  -7.38%  "
$ws.Range("E34").Value = "  -4.93%  "
$ws.Range("D35").Value = "'0.0748"
$ws.Range("E35").Value = "  -5.49%  "
$ws.Range("E36").Value = "  -2.42%  "
$ws.Range("D37").Value = "'17.55"
$ws.Range("E37").Value = "  -1.75%  "
$ws.Range("E38").Value = "  -4.63%  "
$ws.Range("E39").Value = "  -2.48%  "
$ws.Range("E40").Value = "  -8.17%  "
$ws.Range("E41").Value = "  +2.77%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.994.61"
$ws.Range("E43").Value = "  -0.04%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'19.46"
$ws.Range("E44").Value = "  -10.79%  "
$ws.Range("D45").Value = "'0.0286"
$ws.Range("E45").Value = "  -4.46%  "
$ws.Range("E46").Value = "  -7.67%  "
$ws.Range("D47").Value = "'8.74"
$ws.Range("E47").Value = "  -3.77%  "
$ws.Range("D48").Value = "2.722.68"
$ws.Range("E48").Value = "  -2.26%  "
$ws.Range("D49").Value = "'69.59"
$ws.Range("E49").Value = "  -3.81%  "
$ws.Range("E50").Value = "  -4.05%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "'0.178"
$ws.Range("E51").Value = "  -6.91%  "
